# Adds the "Comandos usados en este proyecto" section (with the list of
# npm / project commands) at the end of the document, matching the
# target diff exactly: one leading blank paragraph, a centered bold
# heading, the command lines (each split into the same runs / proofErr
# markers as authored originally by Word's proofing engine), and a
# trailing blank paragraph.

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$d = $word.ActiveDocument

$global:tailRng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$global:tailRng.Collapse(0)

function Add-Para($innerXml) {
    $global:tailRng.InsertParagraphAfter()
    $newRng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
    $xml = "<w:p $wNs>" + $innerXml + "</w:p>"
    $newRng.InsertXML($xml)
    $global:tailRng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
    $global:tailRng.Collapse(0)
}

# blank paragraph
Add-Para ""

# centered, bold heading
Add-Para '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Comandos usados en este proyecto</w:t></w:r>'

# "# npm create vite@latest"
Add-Para '<w:r><w:t xml:space="preserve"># </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>create</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vite@latest</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# "#tema15-entregable"
Add-Para '<w:r><w:t>#</w:t></w:r><w:r><w:t>tema15-entregable</w:t></w:r>'

# "# React"
Add-Para '<w:r><w:t xml:space="preserve"># </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>React</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# "# JavaScript"
Add-Para '<w:r><w:t xml:space="preserve"># </w:t></w:r><w:r><w:t>JavaScript</w:t></w:r>'

# "cd .\tema15-entregable\"
Add-Para '<w:proofErr w:type="gramStart"/><w:r><w:t>cd .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>\tema15-entregable\</w:t></w:r>'

# "npm install"
Add-Para '<w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>install</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# "npm install @mui/material @emotion/react @emotion/styled"
Add-Para '<w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>install</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> @mui/material @emotion/react @emotion/styled</w:t></w:r>'

# "npm install axios"
Add-Para '<w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>install</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>axios</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# "npm install react-router-dom"
Add-Para '<w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>install</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>react</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>router</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-dom</w:t></w:r>'

# trailing blank paragraph
Add-Para ""
